$d = $word.ActiveDocument
$s = $d.Styles.Add("TestSemi", 1)
try {
  $s.SemiHidden = $true
  Write-Output "set semihidden ok"
} catch {
  Write-Output "error: $_"
}
